# Re-run SGNN to annotate dialog acts following clean up work to the original transcripts.
# Updates DAMSLTag (column I) and DialogAct (column J) values for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{Row=4; I="ba"; J="Appreciation"},
    @{Row=23; I="sd"; J="Statement-non-opinion"},
    @{Row=24; I="ba"; J="Appreciation"},
    @{Row=39; I="sv"; J="Statement-opinion"},
    @{Row=43; I="qy"; J="Yes-No-Question"},
    @{Row=44; I="ba"; J="Appreciation"},
    @{Row=50; I="%"; J="Uninterpretable"},
    @{Row=53; I="sv"; J="Statement-opinion"},
    @{Row=55; I="b"; J="Acknowledge (Backchannel)"},
    @{Row=58; I="sv"; J="Statement-opinion"},
    @{Row=88; I="sd"; J="Statement-non-opinion"},
    @{Row=94; I="aa"; J="Agree/Accept"},
    @{Row=97; I="aa"; J="Agree/Accept"},
    @{Row=99; I="aa"; J="Agree/Accept"},
    @{Row=100; I="sd"; J="Statement-non-opinion"},
    @{Row=102; I="b"; J="Acknowledge (Backchannel)"},
    @{Row=103; I="sv"; J="Statement-opinion"},
    @{Row=107; I="aa"; J="Agree/Accept"},
    @{Row=114; I="aa"; J="Agree/Accept"},
    @{Row=115; I="aa"; J="Agree/Accept"},
    @{Row=118; I="sv"; J="Statement-opinion"},
    @{Row=122; I="%"; J="Uninterpretable"},
    @{Row=123; I="aa"; J="Agree/Accept"},
    @{Row=125; I="aa"; J="Agree/Accept"},
    @{Row=128; I="sd"; J="Statement-non-opinion"},
    @{Row=136; I="sd"; J="Statement-non-opinion"},
    @{Row=137; I="sd"; J="Statement-non-opinion"},
    @{Row=140; I="sd"; J="Statement-non-opinion"},
    @{Row=148; I="sv"; J="Statement-opinion"},
    @{Row=150; I="sd"; J="Statement-non-opinion"},
    @{Row=155; I="ba"; J="Appreciation"},
    @{Row=158; I="sv"; J="Statement-opinion"},
    @{Row=160; I="ba"; J="Appreciation"},
    @{Row=165; I="sd"; J="Statement-non-opinion"},
    @{Row=174; I="sd"; J="Statement-non-opinion"},
    @{Row=190; I="sv"; J="Statement-opinion"},
    @{Row=191; I="ba"; J="Appreciation"},
    @{Row=193; I="aa"; J="Agree/Accept"},
    @{Row=197; I="sv"; J="Statement-opinion"},
    @{Row=203; I="b"; J="Acknowledge (Backchannel)"},
    @{Row=208; I="b"; J="Acknowledge (Backchannel)"},
    @{Row=210; I="ba"; J="Appreciation"},
    @{Row=213; I="ba"; J="Appreciation"},
    @{Row=217; I="ba"; J="Appreciation"},
    @{Row=222; I="ba"; J="Appreciation"},
    @{Row=229; I="aa"; J="Agree/Accept"},
    @{Row=237; I="ba"; J="Appreciation"},
    @{Row=252; I="b"; J="Acknowledge (Backchannel)"},
    @{Row=257; I="ba"; J="Appreciation"},
    @{Row=271; I="b"; J="Acknowledge (Backchannel)"},
    @{Row=276; I="aa"; J="Agree/Accept"},
    @{Row=288; I="sd"; J="Statement-non-opinion"},
    @{Row=289; I="b"; J="Acknowledge (Backchannel)"},
    @{Row=293; I="aa"; J="Agree/Accept"},
    @{Row=303; I="%"; J="Uninterpretable"},
    @{Row=306; I="ba"; J="Appreciation"},
    @{Row=310; I="ba"; J="Appreciation"},
    @{Row=318; I="sd"; J="Statement-non-opinion"},
    @{Row=319; I="b"; J="Acknowledge (Backchannel)"},
    @{Row=327; I="b"; J="Acknowledge (Backchannel)"},
    @{Row=339; I="aa"; J="Agree/Accept"},
    @{Row=344; I="ba"; J="Appreciation"},
    @{Row=346; I="ba"; J="Appreciation"},
    @{Row=401; I="sv"; J="Statement-opinion"},
    @{Row=406; I="b"; J="Acknowledge (Backchannel)"},
    @{Row=417; I="sv"; J="Statement-opinion"},
    @{Row=424; I="ba"; J="Appreciation"},
    @{Row=425; I="b"; J="Acknowledge (Backchannel)"},
    @{Row=427; I="aa"; J="Agree/Accept"},
    @{Row=432; I="%"; J="Uninterpretable"},
    @{Row=436; I="ba"; J="Appreciation"},
    @{Row=450; I="ba"; J="Appreciation"},
    @{Row=451; I="aa"; J="Agree/Accept"},
    @{Row=453; I="b"; J="Acknowledge (Backchannel)"},
    @{Row=463; I="b"; J="Acknowledge (Backchannel)"},
    @{Row=471; I="ba"; J="Appreciation"},
    @{Row=473; I="aa"; J="Agree/Accept"},
    @{Row=477; I="%"; J="Uninterpretable"},
    @{Row=485; I="sd"; J="Statement-non-opinion"},
    @{Row=489; I="ba"; J="Appreciation"},
    @{Row=493; I="ba"; J="Appreciation"},
    @{Row=494; I="b"; J="Acknowledge (Backchannel)"},
    @{Row=503; I="sd"; J="Statement-non-opinion"},
    @{Row=504; I="aa"; J="Agree/Accept"},
    @{Row=513; I="sd"; J="Statement-non-opinion"},
    @{Row=520; I="%"; J="Uninterpretable"},
    @{Row=525; I="sd"; J="Statement-non-opinion"},
    @{Row=526; I="sd"; J="Statement-non-opinion"}
)

foreach ($change in $changes) {
    $ws.Cells.Item($change.Row, 9).Value = $change.I
    $ws.Cells.Item($change.Row, 10).Value = $change.J
}
